$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("U1").Value = "inka-mac"
$ws.Range("U3").Value = "/Users/inka/Documents/PrincetonJP/RawData "
$ws.Range("U4").Value = "/Users/inka/Documents/PrincetonJP/FishToolbox "
$ws.Range("U2").Value = "inkab"
$ws.Range("U5").Value = "/Users/inka/Dropbox/LivemRNADatabase"
$ws.Range("U11").Value = "/Users/inka/Documents/PrincetonJP/mRNADynamics"
$ws.Range("U11").Select() | Out-Null
